# Applies the cryptos price/volume update described in the commit
# "Updated cryptos list on Wed Dec 27 03:47:30 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D-column (Price) cells as Text to preserve exact string formatting
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.402.96'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.220.60'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '110.21'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.58'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.621'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.592'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.02'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0906'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.54'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.78'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.546.45'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.224.07'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.293.91'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.02'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.67'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.37'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.35'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '229.04'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.83'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.30'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.00'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.72'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.12'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.71'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0870'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.55'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.89'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.10'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0365'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.103'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.83'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.228'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.15'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.29'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.28'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.70'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.26'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.45'

# Set other (B, C, E) columns
$ws.Range("E2").Value = '  -2.82%  '
$ws.Range("E3").Value = '  -2.16%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("E5").Value = '  -7.79%  '
$ws.Range("E6").Value = '  +7.47%  '
$ws.Range("E7").Value = '  -3.38%  '
$ws.Range("E8").Value = '  -0.41%  '
$ws.Range("E9").Value = '  -4.62%  '
$ws.Range("E10").Value = '  -9.34%  '
$ws.Range("E11").Value = '  -4.04%  '
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("E13").Value = '  -8.67%  '
$ws.Range("E14").Value = '  +8.66%  '
$ws.Range("E15").Value = '  -3.20%  '
$ws.Range("E16").Value = '  -6.45%  '
$ws.Range("E17").Value = '  -2.51%  '
$ws.Range("E18").Value = '  -2.24%  '
$ws.Range("E19").Value = '  -2.86%  '
$ws.Range("E20").Value = '  -5.30%  '
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("E23").Value = '  +11.30%  '
$ws.Range("E24").Value = '  -1.51%  '
$ws.Range("E25").Value = '  -2.37%  '
$ws.Range("E26").Value = '  -8.02%  '
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("E28").Value = '  -7.90%  '
$ws.Range("E29").Value = '  -2.66%  '
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E31").Value = '  -11.99%  '
$ws.Range("B32").Value = 'WEMIXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E32").Value = '  -7.05%  '
$ws.Range("E33").Value = '  -3.83%  '
$ws.Range("E34").Value = '  -5.17%  '
$ws.Range("E35").Value = '  -3.01%  '
$ws.Range("E36").Value = '  +4.77%  '
$ws.Range("E37").Value = '  -4.11%  '
$ws.Range("E38").Value = '  -4.78%  '
$ws.Range("E39").Value = '  -3.79%  '
$ws.Range("E40").Value = '  -5.77%  '
$ws.Range("E41").Value = '  +1.17%  '
$ws.Range("E42").Value = '  -8.08%  '
$ws.Range("E43").Value = '  -4.88%  '
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("E45").Value = '  -11.13%  '
$ws.Range("E46").Value = '  -6.39%  '
$ws.Range("E47").Value = '  -7.48%  '
$ws.Range("E48").Value = '  +5.31%  '
$ws.Range("E49").Value = '  -1.92%  '
$ws.Range("E50").Value = '  -2.28%  '
$ws.Range("E51").Value = '  -2.85%  '
